$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build out new column D by copying the formatting/value from column C,
# then overwrite D3 with the new "key3" JSON payload.
$ws.Range("C1").Copy($ws.Range("D1"))
$ws.Range("C2").Copy($ws.Range("D2"))
$ws.Range("D2").ClearContents()
$ws.Range("C3").Copy($ws.Range("D3"))
$ws.Range("D3").Value2 = '{"target":{"userName":"hugang","caseId":"hugangのテスト"},"value":"key3"}'

# Match column width/bestfit formatting of column C
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth

# Update the active selection to D8
$ws.Range("D8").Select()
